$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,3
$arr[0,0] = -0.6933320760726929
$arr[0,1] = 1.070083141326904
$arr[0,2] = -0.4738787114620209
$arr[1,0] = -0.694248378276825
$arr[1,1] = 0.7021896243095398
$arr[1,2] = -0.1614211350679397
$arr[2,0] = 0.1959350258111953
$arr[2,1] = 0.2964223623275757
$arr[2,2] = -0.2142609804868698
$arr[3,0] = 0.4120286107063293
$arr[3,1] = -0.2721404731273651
$arr[3,2] = -0.3208569586277008
$arr[4,0] = 0.1794416606426239
$arr[4,1] = 0.0216857157647609
$arr[4,2] = -0.3880521357059479
$arr[5,0] = -0.1267545372247696
$arr[5,1] = -0.0375682115554809
$arr[5,2] = -0.1798998117446899
$arr[6,0] = -0.0123700210824608
$arr[6,1] = 0.0419969856739044
$arr[6,2] = 0.271224170923233
$arr[7,0] = -0.0514653958380222
$arr[7,1] = -0.052381694316864
$arr[7,2] = 0.3119994103908539
$arr[8,0] = -0.1554652005434036
$arr[8,1] = -0.0441350154578685
$arr[8,2] = -0.0074830991216003
$arr[9,0] = -0.2000583708286285
$arr[9,1] = -0.1212567538022995
$arr[9,2] = -0.0207694191485643
$arr[10,0] = -0.1815796941518783
$arr[10,1] = -0.0572686158120632
$arr[10,2] = 0.0864374339580535
$arr[11,0] = -0.0739146918058395
$arr[11,1] = -0.1140790879726409
$arr[11,2] = 0.1067487001419067
$arr[12,0] = -0.0395535230636596
$arr[12,1] = -0.0899499058723449
$arr[12,2] = -0.0404698215425014
$arr[13,0] = -0.0148134818300604
$arr[13,1] = 0.1036943718791008
$arr[13,2] = -0.1157589629292488
$arr[14,0] = 0.5971207618713379
$arr[14,1] = 1.289536476135254
$arr[14,2] = -0.3637702465057373
$arr[15,0] = 1.519069194793701
$arr[15,1] = -0.4518875777721405
$arr[15,2] = -0.6734789609909058
$arr[16,0] = 0.2113593816757202
$arr[16,1] = -0.3769038617610931
$arr[16,2] = 0.4825835525989532
$arr[17,0] = 0.2393064647912979
$arr[17,1] = -0.879187822341919
$arr[17,2] = -0.1872301995754242
$arr[18,0] = -0.1922698318958282
$arr[18,1] = -0.9285151958465576
$arr[18,2] = 0.8594874143600464
$arr[19,0] = -3.570354700088501
$arr[19,1] = -0.7802276611328125
$arr[19,2] = -4.989242076873779
$arr[20,0] = -1.221577763557434
$arr[20,1] = 2.375196695327759
$arr[20,2] = -2.503631114959717
$arr[21,0] = 2.165364503860474
$arr[21,1] = 0.5566509366035461
$arr[21,2] = -0.4453207552433014
$arr[22,0] = 0.2727513313293457
$arr[22,1] = 0.5925393104553223
$arr[22,2] = 0.4948008358478546
$arr[23,0] = -0.2121229618787765
$arr[23,1] = 1.80510675907135
$arr[23,2] = 1.96942949295044
$arr[24,0] = 0.2755002379417419
$arr[24,1] = 1.588096976280212
$arr[24,2] = 2.037540912628174
$arr[25,0] = 0.8869763612747192
$arr[25,1] = 0.8231409192085266
$arr[25,2] = 1.362993121147156
$arr[26,0] = 0.2379320114850998
$arr[26,1] = -0.7533495426177979
$arr[26,2] = 0.1786780804395675
$arr[27,0] = 0.6151412725448608
$arr[27,1] = 1.230893492698669
$arr[27,2] = -0.3686571717262268
$arr[28,0] = 0.1403462886810302
$arr[28,1] = 0.7915286421775818
$arr[28,2] = 0.00137444678694
$arr[29,0] = -0.3019201457500458
$arr[29,1] = 0.041233405470848
$arr[29,2] = -0.0345138870179653

$ws.Range("A2:C31").Value = $arr
